# Apply the "working example for pressure" settings update to the
# gpr_settings "tracks" table.
#
# Row 2 already carries the "example" defaults for several settings
# columns (prob_map_thr, shift_k, kernel_adjust, prob_light_w,
# thr_prob_percentile, thr_gs, thr_as, low_speed_fix, and thr_dur).
# This change copies those same default values down into rows 3-19
# (which previously left those columns blank), and resets thr_dur on
# row 2 itself from 12 back to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tracks")

# thr_dur on row 2 goes from 12 -> 0
$ws.Range("D2").Value = 0

$lastRow = 19

# thr_dur (D), prob_map_thr (N), shift_k (O), kernel_adjust (P),
# prob_light_w (Y), thr_prob_percentile (Z), thr_gs (AA), thr_as (AB),
# low_speed_fix (AC) for rows 3..19 get filled in with the same
# working defaults used on row 2.
$ws.Range("D3:D$lastRow").Value = 0
$ws.Range("N3:N$lastRow").Value = 0.9
$ws.Range("O3:O$lastRow").Value = 0
$ws.Range("P3:P$lastRow").Value = 1.4
$ws.Range("Y3:Y$lastRow").Value = 0.1
$ws.Range("Z3:Z$lastRow").Value = 0.9
$ws.Range("AA3:AA$lastRow").Value = 120
$ws.Range("AB3:AB$lastRow").Value = 100
$ws.Range("AC3:AC$lastRow").Value = 15

# Update the active selection to match the saved workbook state.
$ws.Range("E7").Select()
